# Discard invalid ipv4 ips, inswtich
$wb = $excel.ActiveWorkbook

# Remove Sheet2 entirely
$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Delete()

$ws1 = $wb.Worksheets.Item("Sheet1")

# Add new labels in column E for existing rows
$ws1.Range("E1").Value = "local test"
$ws1.Range("E2").Value = "one client, one switch"

# Add new data row 4 (A4 value + E4 label)
$ws1.Range("A4").Value = 0.67057293853488098
$ws1.Range("E4").Value = "two clients, one switch"

# Add new data row 5 (A5 value only)
$ws1.Range("A5").Value = 0.65625

# Update selection / active cell and make Sheet1 the selected tab
[void]$ws1.Range("E4").Select()
[void]$ws1.Activate()
